# Updates cached market/profit figures on the FFXIV leve-profit workbook's
# per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). These columns
# (H..N = currentAveragePrice.. LeveProfitHQ) hold static numbers refreshed
# by a scheduled market-data pull; no formulas are involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 420.27777
$ws.Cells.Item(28, 9).Value = 485.35715
$ws.Cells.Item(28, 10).Value = 192.5
$ws.Cells.Item(28, 11).Value = 485.35715
$ws.Cells.Item(28, 12).Value = 192.5
$ws.Cells.Item(28, 13).Value = -0.3571499999999901
$ws.Cells.Item(28, 14).Value = -1162.5
# Row 41
$ws.Cells.Item(41, 8).Value = 1492.2727
$ws.Cells.Item(41, 9).Value = 1149.125
$ws.Cells.Item(41, 10).Value = 1815.2354
$ws.Cells.Item(41, 11).Value = 1149.125
$ws.Cells.Item(41, 12).Value = 1815.2354
$ws.Cells.Item(41, 13).Value = -709.125
$ws.Cells.Item(41, 14).Value = -2695.2354
# Row 51
$ws.Cells.Item(51, 8).Value = 3977.9375
$ws.Cells.Item(51, 9).Value = 3199
$ws.Cells.Item(51, 11).Value = 3199
$ws.Cells.Item(51, 13).Value = -2715
# Row 62
$ws.Cells.Item(62, 8).Value = 638964.25
$ws.Cells.Item(62, 10).Value = 840242.5
$ws.Cells.Item(62, 12).Value = 840242.5
$ws.Cells.Item(62, 14).Value = -841490.5
# Row 65
$ws.Cells.Item(65, 8).Value = 638964.25
$ws.Cells.Item(65, 10).Value = 840242.5
$ws.Cells.Item(65, 12).Value = 4201212.5
$ws.Cells.Item(65, 14).Value = -4207452.5
# Row 70
$ws.Cells.Item(70, 8).Value = 1133
$ws.Cells.Item(70, 9).Value = 500
$ws.Cells.Item(70, 10).Value = 1259.6
$ws.Cells.Item(70, 11).Value = 1500
$ws.Cells.Item(70, 12).Value = 3778.8
$ws.Cells.Item(70, 13).Value = -1230
$ws.Cells.Item(70, 14).Value = -4318.799999999999
# Row 73
$ws.Cells.Item(73, 8).Value = 1133
$ws.Cells.Item(73, 9).Value = 500
$ws.Cells.Item(73, 10).Value = 1259.6
$ws.Cells.Item(73, 11).Value = 1500
$ws.Cells.Item(73, 12).Value = 3778.8
$ws.Cells.Item(73, 13).Value = -564
$ws.Cells.Item(73, 14).Value = -5650.799999999999
# Row 92
$ws.Cells.Item(92, 8).Value = 399.5
$ws.Cells.Item(92, 9).Value = 439.5
$ws.Cells.Item(92, 11).Value = 439.5
$ws.Cells.Item(92, 13).Value = 808.5
# Row 94
$ws.Cells.Item(94, 8).Value = 6586.4443
$ws.Cells.Item(94, 9).Value = 6586.4443
$ws.Cells.Item(94, 11).Value = 6586.4443
$ws.Cells.Item(94, 13).Value = -6135.4443
# Row 96
$ws.Cells.Item(96, 8).Value = 640.7646999999999
$ws.Cells.Item(96, 9).Value = 646.3333
$ws.Cells.Item(96, 11).Value = 1938.9999
$ws.Cells.Item(96, 13).Value = -565.9999
# Row 98
$ws.Cells.Item(98, 8).Value = 1250.6666
$ws.Cells.Item(98, 9).Value = 1250.6666
$ws.Cells.Item(98, 11).Value = 1250.6666
$ws.Cells.Item(98, 13).Value = 247.3334
# Row 99
$ws.Cells.Item(99, 8).Value = 354.08334
$ws.Cells.Item(99, 9).Value = 244.9
$ws.Cells.Item(99, 10).Value = 900
$ws.Cells.Item(99, 11).Value = 734.7
$ws.Cells.Item(99, 12).Value = 2700
$ws.Cells.Item(99, 13).Value = 763.3
$ws.Cells.Item(99, 14).Value = -5696
# Row 107
$ws.Cells.Item(107, 8).Value = 57996.723
$ws.Cells.Item(107, 9).Value = 86402.75
$ws.Cells.Item(107, 10).Value = 1184.6666
$ws.Cells.Item(107, 11).Value = 86402.75
$ws.Cells.Item(107, 12).Value = 1184.6666
$ws.Cells.Item(107, 13).Value = -84482.75
$ws.Cells.Item(107, 14).Value = -5024.6666
# Row 115
$ws.Cells.Item(115, 8).Value = 650.7222
$ws.Cells.Item(115, 9).Value = 703.9375
$ws.Cells.Item(115, 10).Value = 225
$ws.Cells.Item(115, 11).Value = 2111.8125
$ws.Cells.Item(115, 12).Value = 675
$ws.Cells.Item(115, 13).Value = -544.8125
$ws.Cells.Item(115, 14).Value = -3809
# Row 116
$ws.Cells.Item(116, 8).Value = 10224.375
$ws.Cells.Item(116, 9).Value = 11999.25
$ws.Cells.Item(116, 11).Value = 11999.25
$ws.Cells.Item(116, 13).Value = -8557.25
# Row 118
$ws.Cells.Item(118, 8).Value = 775.1
$ws.Cells.Item(118, 9).Value = 710.6316
$ws.Cells.Item(118, 11).Value = 2131.8948
$ws.Cells.Item(118, 13).Value = -474.8948
# Row 122
$ws.Cells.Item(122, 8).Value = 1250.6666
$ws.Cells.Item(122, 9).Value = 1250.6666
$ws.Cells.Item(122, 11).Value = 3751.9998
$ws.Cells.Item(122, 13).Value = -1301.9998
# Row 132
$ws.Cells.Item(132, 8).Value = 2425.7407
$ws.Cells.Item(132, 9).Value = 1395.9565
$ws.Cells.Item(132, 10).Value = 8347
$ws.Cells.Item(132, 11).Value = 4187.8695
$ws.Cells.Item(132, 12).Value = 25041
$ws.Cells.Item(132, 13).Value = -1657.8695
$ws.Cells.Item(132, 14).Value = -30101
# Row 135
$ws.Cells.Item(135, 8).Value = 752.7917
$ws.Cells.Item(135, 10).Value = 1613.2858
$ws.Cells.Item(135, 12).Value = 14519.5722
$ws.Cells.Item(135, 14).Value = -19589.5722
# Row 137
$ws.Cells.Item(137, 8).Value = 333681.1
$ws.Cells.Item(137, 10).Value = 715911.0600000001
$ws.Cells.Item(137, 12).Value = 2147733.18
$ws.Cells.Item(137, 14).Value = -2152833.18
# Row 138
$ws.Cells.Item(138, 8).Value = 2546.9827
$ws.Cells.Item(138, 9).Value = 2003.3636
$ws.Cells.Item(138, 10).Value = 2674.2126
$ws.Cells.Item(138, 11).Value = 6010.0908
$ws.Cells.Item(138, 12).Value = 8022.6378
$ws.Cells.Item(138, 13).Value = -870.0907999999999
$ws.Cells.Item(138, 14).Value = -18302.6378

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6631.344
$ws.Cells.Item(32, 9).Value = 3825.6375
$ws.Cells.Item(32, 10).Value = 23897.23
$ws.Cells.Item(32, 11).Value = 3825.6375
$ws.Cells.Item(32, 12).Value = 23897.23
$ws.Cells.Item(32, 13).Value = -3538.6375
$ws.Cells.Item(32, 14).Value = -24471.23
# Row 44
$ws.Cells.Item(44, 8).Value = 122288
$ws.Cells.Item(44, 10).Value = 122288
$ws.Cells.Item(44, 12).Value = 122288
$ws.Cells.Item(44, 14).Value = -123264
# Row 45
$ws.Cells.Item(45, 8).Value = 14225
$ws.Cells.Item(45, 9).Value = 15828.571
$ws.Cells.Item(45, 10).Value = 3000
$ws.Cells.Item(45, 11).Value = 15828.571
$ws.Cells.Item(45, 12).Value = 3000
$ws.Cells.Item(45, 13).Value = -15451.571
$ws.Cells.Item(45, 14).Value = -3754
# Row 74
$ws.Cells.Item(74, 8).Value = 3354.158
$ws.Cells.Item(74, 9).Value = 2016.2142
$ws.Cells.Item(74, 10).Value = 7100.4
$ws.Cells.Item(74, 11).Value = 2016.2142
$ws.Cells.Item(74, 12).Value = 7100.4
$ws.Cells.Item(74, 13).Value = -1142.2142
$ws.Cells.Item(74, 14).Value = -8848.4
# Row 77
$ws.Cells.Item(77, 8).Value = 3354.158
$ws.Cells.Item(77, 9).Value = 2016.2142
$ws.Cells.Item(77, 10).Value = 7100.4
$ws.Cells.Item(77, 11).Value = 10081.071
$ws.Cells.Item(77, 12).Value = 35502
$ws.Cells.Item(77, 13).Value = -5713.071
$ws.Cells.Item(77, 14).Value = -44238
# Row 86
$ws.Cells.Item(86, 8).Value = 19999
$ws.Cells.Item(86, 9).Value = 19999
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 19999
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -18813
$ws.Cells.Item(86, 14).ClearContents()
# Row 89
$ws.Cells.Item(89, 8).Value = 19999
$ws.Cells.Item(89, 9).Value = 19999
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 59997
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -54069
$ws.Cells.Item(89, 14).ClearContents()
# Row 97
$ws.Cells.Item(97, 8).Value = 3538.182
$ws.Cells.Item(97, 9).Value = 878.8333
$ws.Cells.Item(97, 11).Value = 878.8333
$ws.Cells.Item(97, 13).Value = -382.8333
# Row 110
$ws.Cells.Item(110, 8).Value = 1075
$ws.Cells.Item(110, 9).Value = 1066.375
$ws.Cells.Item(110, 11).Value = 1066.375
$ws.Cells.Item(110, 13).Value = 978.625
# Row 132
$ws.Cells.Item(132, 8).Value = 2633.0303
$ws.Cells.Item(132, 9).Value = 2354.0417
$ws.Cells.Item(132, 10).Value = 3377
$ws.Cells.Item(132, 11).Value = 7062.125100000001
$ws.Cells.Item(132, 12).Value = 10131
$ws.Cells.Item(132, 13).Value = -4532.125100000001
$ws.Cells.Item(132, 14).Value = -15191

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Cells.Item(64, 8).Value = 1189
$ws.Cells.Item(64, 9).Value = 99.5
$ws.Cells.Item(64, 11).Value = 99.5
$ws.Cells.Item(64, 13).Value = 125.5
# Row 67
$ws.Cells.Item(67, 8).Value = 1189
$ws.Cells.Item(67, 9).Value = 99.5
$ws.Cells.Item(67, 11).Value = 99.5
$ws.Cells.Item(67, 13).Value = 680.5
# Row 86
$ws.Cells.Item(86, 8).Value = 1810.1852
$ws.Cells.Item(86, 9).Value = 1431.2941
$ws.Cells.Item(86, 10).Value = 2454.3
$ws.Cells.Item(86, 11).Value = 1431.2941
$ws.Cells.Item(86, 12).Value = 2454.3
$ws.Cells.Item(86, 13).Value = -308.2941000000001
$ws.Cells.Item(86, 14).Value = -4700.3
# Row 89
$ws.Cells.Item(89, 8).Value = 1810.1852
$ws.Cells.Item(89, 9).Value = 1431.2941
$ws.Cells.Item(89, 10).Value = 2454.3
$ws.Cells.Item(89, 11).Value = 7156.4705
$ws.Cells.Item(89, 12).Value = 12271.5
$ws.Cells.Item(89, 13).Value = -1540.4705
$ws.Cells.Item(89, 14).Value = -23503.5
# Row 94
$ws.Cells.Item(94, 8).Value = 1453.0714
$ws.Cells.Item(94, 9).Value = 1407.5555
$ws.Cells.Item(94, 10).Value = 1535
$ws.Cells.Item(94, 11).Value = 1407.5555
$ws.Cells.Item(94, 12).Value = 1535
$ws.Cells.Item(94, 13).Value = -956.5554999999999
$ws.Cells.Item(94, 14).Value = -2437
# Row 99
$ws.Cells.Item(99, 8).Value = 2265.4285
$ws.Cells.Item(99, 9).Value = 1919.4
$ws.Cells.Item(99, 10).Value = 3130.5
$ws.Cells.Item(99, 11).Value = 1919.4
$ws.Cells.Item(99, 12).Value = 3130.5
$ws.Cells.Item(99, 13).Value = -421.4000000000001
$ws.Cells.Item(99, 14).Value = -6126.5
# Row 107
$ws.Cells.Item(107, 8).Value = 1996.919
$ws.Cells.Item(107, 9).Value = 1822.2354
$ws.Cells.Item(107, 10).Value = 3976.6667
$ws.Cells.Item(107, 11).Value = 1822.2354
$ws.Cells.Item(107, 12).Value = 3976.6667
$ws.Cells.Item(107, 13).Value = 97.76459999999997
$ws.Cells.Item(107, 14).Value = -7816.6667
# Row 134
$ws.Cells.Item(134, 8).Value = 1536.2122
$ws.Cells.Item(134, 9).Value = 1167.75
$ws.Cells.Item(134, 11).Value = 3503.25
$ws.Cells.Item(134, 13).Value = -968.25

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 440.75
$ws.Cells.Item(22, 9).Value = 440.75
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 440.75
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -90.75
$ws.Cells.Item(22, 14).ClearContents()
# Row 31
$ws.Cells.Item(31, 8).Value = 3385.0889
$ws.Cells.Item(31, 9).Value = 1996.4166
$ws.Cells.Item(31, 10).Value = 4972.143
$ws.Cells.Item(31, 11).Value = 1996.4166
$ws.Cells.Item(31, 12).Value = 4972.143
$ws.Cells.Item(31, 13).Value = -1701.4166
$ws.Cells.Item(31, 14).Value = -5562.143
# Row 34
$ws.Cells.Item(34, 8).Value = 3385.0889
$ws.Cells.Item(34, 9).Value = 1996.4166
$ws.Cells.Item(34, 10).Value = 4972.143
$ws.Cells.Item(34, 11).Value = 1996.4166
$ws.Cells.Item(34, 12).Value = 4972.143
$ws.Cells.Item(34, 13).Value = -1794.4166
$ws.Cells.Item(34, 14).Value = -5376.143
# Row 111
$ws.Cells.Item(111, 8).Value = 50699.25
$ws.Cells.Item(111, 10).Value = 50699.25
$ws.Cells.Item(111, 12).Value = 50699.25
$ws.Cells.Item(111, 14).Value = -58879.25
# Row 141
$ws.Cells.Item(141, 8).Value = 124069.93
$ws.Cells.Item(141, 10).Value = 124069.93
$ws.Cells.Item(141, 12).Value = 124069.93
$ws.Cells.Item(141, 14).Value = -134429.93

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Cells.Item(40, 8).Value = 123.4
$ws.Cells.Item(40, 9).Value = 73
$ws.Cells.Item(40, 11).Value = 292
$ws.Cells.Item(40, 13).Value = -223
# Row 50
$ws.Cells.Item(50, 8).Value = 1087.25
$ws.Cells.Item(50, 9).Value = 1049.6666
$ws.Cells.Item(50, 11).Value = 3148.9998
$ws.Cells.Item(50, 13).Value = -2667.9998
# Row 53
$ws.Cells.Item(53, 8).Value = 1087.25
$ws.Cells.Item(53, 9).Value = 1049.6666
$ws.Cells.Item(53, 11).Value = 3148.9998
$ws.Cells.Item(53, 13).Value = -2667.9998
# Row 61
$ws.Cells.Item(61, 8).Value = 174.92857
$ws.Cells.Item(61, 9).Value = 211.55556
$ws.Cells.Item(61, 10).Value = 109
$ws.Cells.Item(61, 11).Value = 634.66668
$ws.Cells.Item(61, 12).Value = 327
$ws.Cells.Item(61, 13).Value = -419.66668
$ws.Cells.Item(61, 14).Value = -757
# Row 68
$ws.Cells.Item(68, 8).Value = 1433.3334
$ws.Cells.Item(68, 10).Value = 650
$ws.Cells.Item(68, 12).Value = 1950
$ws.Cells.Item(68, 14).Value = -3572
# Row 71
$ws.Cells.Item(71, 8).Value = 1433.3334
$ws.Cells.Item(71, 10).Value = 650
$ws.Cells.Item(71, 12).Value = 5850
$ws.Cells.Item(71, 14).Value = -13962
# Row 112
$ws.Cells.Item(112, 8).Value = 5167.3477
$ws.Cells.Item(112, 9).Value = 850
$ws.Cells.Item(112, 11).Value = 2550
$ws.Cells.Item(112, 13).Value = -1442
# Row 117
$ws.Cells.Item(117, 8).Value = 3999.6667
$ws.Cells.Item(117, 9).Value = 1999
$ws.Cells.Item(117, 10).Value = 5000
$ws.Cells.Item(117, 11).Value = 5997
$ws.Cells.Item(117, 12).Value = 15000
$ws.Cells.Item(117, 13).Value = -2555
$ws.Cells.Item(117, 14).Value = -21884
# Row 120
$ws.Cells.Item(120, 8).Value = 36838.168
$ws.Cells.Item(120, 9).Value = 23676.666
$ws.Cells.Item(120, 11).Value = 71029.99800000001
$ws.Cells.Item(120, 13).Value = -66191.99800000001
# Row 122
$ws.Cells.Item(122, 8).Value = 944
$ws.Cells.Item(122, 9).Value = 983.3333
$ws.Cells.Item(122, 10).Value = 885
$ws.Cells.Item(122, 11).Value = 8849.9997
$ws.Cells.Item(122, 12).Value = 7965
$ws.Cells.Item(122, 13).Value = -6399.9997
$ws.Cells.Item(122, 14).Value = -12865
# Row 132
$ws.Cells.Item(132, 8).Value = 36259.465
$ws.Cells.Item(132, 10).Value = 44695.793
$ws.Cells.Item(132, 12).Value = 402262.137
$ws.Cells.Item(132, 14).Value = -407322.137

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Cells.Item(21, 8).Value = 430366.66
$ws.Cells.Item(21, 9).Value = 1111200
$ws.Cells.Item(21, 10).Value = 21866.666
$ws.Cells.Item(21, 11).Value = 1111200
$ws.Cells.Item(21, 12).Value = 21866.666
$ws.Cells.Item(21, 13).Value = -1111027
$ws.Cells.Item(21, 14).Value = -22212.666
# Row 29
$ws.Cells.Item(29, 8).Value = 5000
$ws.Cells.Item(29, 10).Value = 5000
$ws.Cells.Item(29, 12).Value = 5000
$ws.Cells.Item(29, 14).Value = -5580
# Row 30
$ws.Cells.Item(30, 8).Value = 430366.66
$ws.Cells.Item(30, 9).Value = 1111200
$ws.Cells.Item(30, 10).Value = 21866.666
$ws.Cells.Item(30, 11).Value = 1111200
$ws.Cells.Item(30, 12).Value = 21866.666
$ws.Cells.Item(30, 13).Value = -1111095
$ws.Cells.Item(30, 14).Value = -22076.666
# Row 107
$ws.Cells.Item(107, 8).Value = 765.1177
$ws.Cells.Item(107, 9).Value = 759.0909
$ws.Cells.Item(107, 10).Value = 776.1667
$ws.Cells.Item(107, 11).Value = 759.0909
$ws.Cells.Item(107, 12).Value = 776.1667
$ws.Cells.Item(107, 13).Value = 1160.9091
$ws.Cells.Item(107, 14).Value = -4616.1667
# Row 113
$ws.Cells.Item(113, 8).Value = 2163.875
$ws.Cells.Item(113, 9).Value = 1133.3334
$ws.Cells.Item(113, 10).Value = 2782.2
$ws.Cells.Item(113, 11).Value = 1133.3334
$ws.Cells.Item(113, 12).Value = 2782.2
$ws.Cells.Item(113, 13).Value = 1036.6666
$ws.Cells.Item(113, 14).Value = -7122.2
# Row 122
$ws.Cells.Item(122, 8).Value = 6057.725
$ws.Cells.Item(122, 9).Value = 5379.794
$ws.Cells.Item(122, 11).Value = 16139.382
$ws.Cells.Item(122, 13).Value = -13689.382

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 2064.3333
$ws.Cells.Item(22, 9).Value = 567.25
$ws.Cells.Item(22, 11).Value = 567.25
$ws.Cells.Item(22, 13).Value = -272.25
# Row 26
$ws.Cells.Item(26, 8).Value = 1603.75
$ws.Cells.Item(26, 10).Value = 1805
$ws.Cells.Item(26, 12).Value = 1805
$ws.Cells.Item(26, 14).Value = -2395
# Row 27
$ws.Cells.Item(27, 8).Value = 2064.3333
$ws.Cells.Item(27, 9).Value = 567.25
$ws.Cells.Item(27, 11).Value = 567.25
$ws.Cells.Item(27, 13).Value = -460.25
# Row 31
$ws.Cells.Item(31, 8).Value = 2356.8
$ws.Cells.Item(31, 10).Value = 3730.3333
$ws.Cells.Item(31, 12).Value = 3730.3333
$ws.Cells.Item(31, 14).Value = -4226.3333
# Row 93
$ws.Cells.Item(93, 8).Value = 3098.8
$ws.Cells.Item(93, 9).Value = 3347.8
$ws.Cells.Item(93, 10).Value = 2974.3
$ws.Cells.Item(93, 11).Value = 3347.8
$ws.Cells.Item(93, 12).Value = 2974.3
$ws.Cells.Item(93, 13).Value = -2099.8
$ws.Cells.Item(93, 14).Value = -5470.3
# Row 132
$ws.Cells.Item(132, 8).Value = 3982.3333
$ws.Cells.Item(132, 9).Value = 3260.2
$ws.Cells.Item(132, 10).Value = 4498.143
$ws.Cells.Item(132, 11).Value = 9780.599999999999
$ws.Cells.Item(132, 12).Value = 13494.429
$ws.Cells.Item(132, 13).Value = -7250.599999999999
$ws.Cells.Item(132, 14).Value = -18554.429
# Row 136
$ws.Cells.Item(136, 8).Value = 6946.5835
$ws.Cells.Item(136, 9).Value = 7162.3335
$ws.Cells.Item(136, 10).Value = 6299.3335
$ws.Cells.Item(136, 11).Value = 21487.0005
$ws.Cells.Item(136, 12).Value = 18898.0005
$ws.Cells.Item(136, 13).Value = -18937.0005
$ws.Cells.Item(136, 14).Value = -23998.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Cells.Item(31, 8).Value = 50025000
$ws.Cells.Item(31, 10).Value = 50001
$ws.Cells.Item(31, 12).Value = 50001
$ws.Cells.Item(31, 14).Value = -50697
# Row 54
$ws.Cells.Item(54, 8).Value = 38479.8
$ws.Cells.Item(54, 10).Value = 37666.332
$ws.Cells.Item(54, 12).Value = 37666.332
$ws.Cells.Item(54, 14).Value = -38706.332
# Row 95
$ws.Cells.Item(95, 8).Value = 23344
$ws.Cells.Item(95, 10).Value = 23344
$ws.Cells.Item(95, 12).Value = 23344
$ws.Cells.Item(95, 14).Value = -28836
# Row 105
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).ClearContents()
# Row 123
$ws.Cells.Item(123, 8).Value = 79000
$ws.Cells.Item(123, 10).Value = 79000
$ws.Cells.Item(123, 12).Value = 79000
$ws.Cells.Item(123, 14).Value = -88800
# Row 132
$ws.Cells.Item(132, 8).Value = 1452045.2
$ws.Cells.Item(132, 9).Value = 2499.2856
$ws.Cells.Item(132, 10).Value = 4834319.5
$ws.Cells.Item(132, 11).Value = 7497.8568
$ws.Cells.Item(132, 12).Value = 14502958.5
$ws.Cells.Item(132, 13).Value = -4967.8568
$ws.Cells.Item(132, 14).Value = -14508018.5

